# Applies the cryptos.xlsx price/volume update described in the commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.013.88'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.562.45'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("D4").Value = '''1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''207.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("E6").Value = '  +0.97%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '''22.16'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.04%  '
$ws.Range("D9").Value = '''0.249'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("D10").Value = '''0.0597'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.79%  '
$ws.Range("D11").Value = '''0.0860'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("D12").Value = '1.782.26'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").Value = '1.539.92'
$ws.Range("E13").Value = '  -0.69%  '
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").Value = '''62.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("D17").Value = '27.003.15'
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("E18").Value = '  +2.76%  '
$ws.Range("D19").Value = '''217.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("D20").Value = '''7.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.10%  '
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("E22").Value = '  +1.49%  '
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("E24").Value = '  -2.67%  '
$ws.Range("D25").Value = '''153.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  +1.18%  '
$ws.Range("E28").Value = '  +1.51%  '
$ws.Range("D29").Value = '''1.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").Value = '''0.0469'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("E31").Value = '  +1.73%  '
$ws.Range("E32").Value = '  +0.41%  '
$ws.Range("D33").Value = '1.423.61'
$ws.Range("E33").Value = '  +0.42%  '
$ws.Range("E34").Value = '  +3.58%  '
$ws.Range("E35").Value = '  +2.77%  '
$ws.Range("E36").Value = '  +8.72%  '
$ws.Range("E37").Value = '  +1.38%  '
$ws.Range("E38").Value = '  +0.68%  '
$ws.Range("D39").Value = '''0.531'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.84%  '
$ws.Range("D40").Value = '''0.810'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.40%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").Value = '''2.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.91%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.67'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.65%  '
$ws.Range("D44").Value = '''1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.70%  '
$ws.Range("D45").Value = '''65.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.04%  '
$ws.Range("D46").Value = '''1.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").Value = '1.698.03'
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("D48").Value = '''87.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.74%  '
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("D50").Value = '0.0₇0993'
$ws.Range("E50").Value = '  -1.24%  '
$ws.Range("E51").Value = '  -0.22%  '
